# "Group Reflection OneDrive ReSync"
#
# 1. The "At least one thing that you have learned about groups?" heading
#    paragraph loses its proofing (grammar) marks and is retyped as two
#    separate (but identically-bold) runs:
#       "At least one thing that you have learned about " + "groups?"
# 2. The section's page size is stamped with an explicit (portrait) orientation.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the "...learned about groups?" heading paragraph and replace
#    its contents (text only - leave the paragraph mark itself alone) with
#    two bold runs, dropping the old proofErr gramStart/gramEnd wrapper.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*learned about groups?*") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $paraRange = $target.Range

    $openXmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/temp.xml" pkg:contentType="application/xml">' +
        '<pkg:xmlData>' +
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>' +
        '<w:t xml:space="preserve">At least one thing that you have learned about </w:t></w:r>' +
        '<w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>' +
        '<w:t>groups?</w:t></w:r>' +
        '</w:p>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $paraRange.InsertXML($openXmlFragment)
}

# ---------------------------------------------------------------------
# 2. Make the page orientation explicit (portrait) on every section.
# ---------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $d.Sections.Item($s).PageSetup.Orientation = 0
}
